# Applies updated crypto price/volume data per GitHub Actions scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''61.734.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.41%  '

$ws.Range("D3").Value = '''2.580.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.90%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '''553.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.81%  '

$ws.Range("D6").Value = '''153.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.96%  '

$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = '''0.596'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.63%  '

$ws.Range("E9").Value = '  -2.54%  '

$ws.Range("E10").Value = '  -2.00%  '

$ws.Range("D11").Value = '''5.48'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.28%  '

$ws.Range("D12").Value = '''0.363'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.40%  '

$ws.Range("D13").Value = '''3.038.50'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.84%  '

$ws.Range("D14").Value = '''25.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.80%  '

$ws.Range("D15").Value = '''61.656.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.37%  '

$ws.Range("D16").Value = '''0.0000143'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.11%  '

$ws.Range("D17").Value = '''2.585.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.79%  '

$ws.Range("D18").Value = '''11.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.77%  '

$ws.Range("D19").Value = '''4.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.66%  '

$ws.Range("D20").Value = '''338.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.22%  '

$ws.Range("D21").Value = '''6.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.32%  '

$ws.Range("D22").Value = '''0.998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("D23").Value = '''0.494'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.75%  '

$ws.Range("D24").Value = '''62.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.07%  '

$ws.Range("E25").Value = '  -0.26%  '

$ws.Range("D26").Value = '''0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.21%  '

$ws.Range("D27").Value = '''8.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.45%  '

$ws.Range("D28").Value = '''0.0₃0836'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.15%  '

$ws.Range("E29").Value = '  -1.57%  '

$ws.Range("D30").Value = '''7.04'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.21%  '

$ws.Range("E31").Value = '  -5.27%  '

$ws.Range("D32").Value = '''160.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.43%  '

$ws.Range("D33").Value = '''0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.04%  '

$ws.Range("D34").Value = '''19.21'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.88%  '

$ws.Range("D35").Value = '''4.67'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.77%  '

$ws.Range("D36").Value = '''1.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.34%  '

$ws.Range("D37").Value = '''1.79'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.10%  '

$ws.Range("D38").Value = '''336.27'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.83%  '

$ws.Range("D39").Value = '''6.04'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.48%  '

$ws.Range("D40").Value = '''0.891'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.85%  '

$ws.Range("D41").Value = '''3.92'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.99%  '

$ws.Range("D42").Value = '''37.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.04%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '''0.998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.02%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '''20.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.62%  '

$ws.Range("D45").Value = '''2.131.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.25%  '

$ws.Range("D46").Value = '''0.606'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.40%  '

$ws.Range("D47").Value = '''10.94'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.97%  '

$ws.Range("D48").Value = '''19.57'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.68%  '

$ws.Range("D49").Value = '''0.0545'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.37%  '

$ws.Range("E50").Value = '  -1.82%  '

$ws.Range("D51").Value = '''0.0239'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.99%  '
